# CORE_holdings.xlsx update:
#  - bump the "as of" date in the confidentiality footnote (A11) from
#    2021-03-26 to 2021-03-29
#  - refresh the Weight (col D) / Percent Change (col E) figures for the
#    six model sleeves + Total row (rows 2-8)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet ships protected, so values can't be written until it's unlocked.
$ws.Unprotect()

# --- A11: disclaimer footnote text -----------------------------------
$ws.Range("A11").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-03-29 for illustrative purposes only and are subject to change."

# Undo the implicit row-height autofit that comes from re-setting a
# multi-line cell, so row 11 keeps using the sheet's default height.
$ws.Rows(11).RowHeight = 15

# --- D2:E8 -> updated Weight / Percent Change values ------------------
$ws.Range("D2").Value = 0.4986802121063496
$ws.Range("E2").Value = -0.002601420234830898

$ws.Range("D3").Value = 0.2418656242766068
$ws.Range("E3").Value = 0.001083926912356814

$ws.Range("D4").Value = 0.09848326987062692
$ws.Range("E4").Value = -0.01886071109978327

$ws.Range("D5").Value = 0.102425501210619
$ws.Range("E5").Value = -0.01760186715938927

$ws.Range("D6").Value = 0.03043480686074554
$ws.Range("E6").Value = -0.02592047128129615

$ws.Range("D7").Value = 0.02811058567505208
$ws.Range("E7").Value = -0.02463361396944164

$ws.Range("D8").Value = 1
$ws.Range("E8").Value = -0.006176806555550751

# Restore sheet protection (matches the workbook's original, shipped
# protected state -- exact legacy password hash cannot be reproduced via
# the object model, which only exposes the modern salted-hash Protect()).
$ws.Protect()
